$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) First paragraph ("PROCURAÇÃO" title): add w:ind (right=-120, firstLine=0)
#    w:ind values are in twentieths of a point (dxa); ParagraphFormat
#    indents are expressed in points, so -120 dxa == -6 pt, 0 dxa == 0 pt.
# -----------------------------------------------------------------------
$titlePar = $d.Paragraphs.Item(1)
$titlePar.Range.ParagraphFormat.RightIndent = -6
$titlePar.Range.ParagraphFormat.FirstLineIndent = 0

# -----------------------------------------------------------------------
# 2) PODERES paragraph: change the trailing ";" to ". " and append the new
#    "_proc_powers_" placeholder text (Sec. Areas & Powers) as its own run.
# -----------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "reserva de poderes;", $true, $false, $false, $false, $false,
    $true, 1, $false, "reserva de poderes. _proc_powers_", 2)

$placeholder = $d.Content.Duplicate
$null = $placeholder.Find.Execute("_proc_powers_")
$placeholder.Font.Name = "Verdana"
$placeholder.Font.Size = 10

# -----------------------------------------------------------------------
# 3) Section page margins: top becomes 0, bottom/left/right become 0.79in.
# -----------------------------------------------------------------------
$d.PageSetup.TopMargin = 0
$d.PageSetup.BottomMargin = 56.88
$d.PageSetup.LeftMargin = 56.88
$d.PageSetup.RightMargin = 56.88
